# Generate Report for Handoff
# Record the latest handoff datetime for the c03cd041 file (row 5)
# on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-26 05:16:45"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-26 05:17:00"
